$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Data rows: six 3-row blocks (DULJINA NIZA label, ZAUZETA MEMORIJA, three
# VRIJEME IZVODENJA trials) plus the average-of-three formula in column E.
# ---------------------------------------------------------------------------

$ws.Range("B3").Value = "10 000 000"
$ws.Range("C3").Value = 4957
$ws.Range("D3").Value = 16219
$ws.Range("D4").Value = 16551
$ws.Range("D5").Value = 15855
$ws.Range("E5").Formula = "=AVERAGE(D3:D5)"

$ws.Range("B6").Value = "4 647 121"
$ws.Range("C6").Value = 2448
$ws.Range("D6").Value = 8129
$ws.Range("D7").Value = 8046
$ws.Range("D8").Value = 7983
$ws.Range("E8").Formula = "=AVERAGE(D6:D8)"

$ws.Range("B9").Value = "1 000 000"
$ws.Range("C9").Value = 541
$ws.Range("D9").Value = 1859
$ws.Range("D10").Value = 1762
$ws.Range("D11").Value = 1869
$ws.Range("E11").Formula = "=AVERAGE(D9:D11)"

$ws.Range("B12").Value = "100 000"
$ws.Range("C12").Value = 61
$ws.Range("D12").Value = 177
$ws.Range("D13").Value = 168
$ws.Range("D14").Value = 167
$ws.Range("E14").Formula = "=AVERAGE(D12:D14)"

$ws.Range("B15").Value = "10 000"
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = 18
$ws.Range("D16").Value = 11
$ws.Range("D17").Value = 16
$ws.Range("E17").Formula = "=AVERAGE(D15:D17)"

$ws.Range("B18").Value = "1 000"
$ws.Range("C18").Value = 5.5
$ws.Range("D18").Value = 1
$ws.Range("D19").Value = 1
$ws.Range("D20").Value = 1
$ws.Range("E20").Formula = "=AVERAGE(D18:D20)"

$ws.Range("B21").Value = 100
$ws.Range("C21").Value = 0.4
$ws.Range("D21").Value = 0.1
$ws.Range("D22").Value = 0.1
$ws.Range("D23").Value = 0.1
$ws.Range("E23").Formula = "=AVERAGE(D21:D23)"

# New "AVG (ms)" column header.
$ws.Range("E2").Value = "AVG (ms)"

# ---------------------------------------------------------------------------
# Formatting: B/C columns centered (both axes), D column left-aligned.
# Build the exact style once on a single cell, then copy/paste the format
# onto the rest of each block so no stray intermediate style gets minted.
# ---------------------------------------------------------------------------

$ws.Range("B3").HorizontalAlignment = -4108
$ws.Range("B3").VerticalAlignment = -4108
$ws.Range("B3").Copy()
$ws.Range("B3:C23").PasteSpecial(-4122)

$ws.Range("D3").HorizontalAlignment = -4131
$ws.Range("D3").Copy()
$ws.Range("D3:D23").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Merge the B and C cells within each 3-row block.
# ---------------------------------------------------------------------------

$ws.Range("B3:B5").Merge()
$ws.Range("C3:C5").Merge()
$ws.Range("B6:B8").Merge()
$ws.Range("C6:C8").Merge()
$ws.Range("B9:B11").Merge()
$ws.Range("C9:C11").Merge()
$ws.Range("B12:B14").Merge()
$ws.Range("C12:C14").Merge()
$ws.Range("B15:B17").Merge()
$ws.Range("C15:C17").Merge()
$ws.Range("B18:B20").Merge()
$ws.Range("C18:C20").Merge()
$ws.Range("B21:B23").Merge()
$ws.Range("C21:C23").Merge()

# ---------------------------------------------------------------------------
# Column E width to match the new data.
# ---------------------------------------------------------------------------

$ws.Columns("E").ColumnWidth = 24.42

# ---------------------------------------------------------------------------
# Selection cursor, matching where the author's cursor ended up.
# ---------------------------------------------------------------------------

$ws.Range("D32").Select()
